$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = '[Chaolin%Huang%NULL%14,                Yeming%Wang%NULL%16,                Xingwang%Li%NULL%7,                Lili%Ren%NULL%6,                Jianping%Zhao%NULL%13,                Yi%Hu%NULL%14,                Li%Zhang%NULL%13,                Guohui%Fan%NULL%16,                Jiuyang%Xu%NULL%16,                Xiaoying%Gu%NULL%16,                Zhenshun%Cheng%NULL%9,                Ting%Yu%NULL%27,                Jiaan%Xia%NULL%7,                Yuan%Wei%NULL%21,                Wenjuan%Wu%NULL%6,                Xuelei%Xie%NULL%7,                Wen%Yin%NULL%10,                Hui%Li%NULL%18,                Min%Liu%NULL%6,                Yan%Xiao%NULL%7,                Hong%Gao%NULL%7,                Li%Guo%NULL%8,                Jungang%Xie%NULL%7,                Guangfa%Wang%NULL%6,                Rongmeng%Jiang%NULL%6,                Zhancheng%Gao%NULL%6,                Qi%Jin%NULL%7,                Jianwei%Wang%wangjw28@163.com%6,                Bin%Cao%caobin_ben@163.com%16]'
$ws.Range("E3").Value = '[Jasper Fuk-Woo%Chan%NULL%3,                Shuofeng%Yuan%NULL%2,                Kin-Hang%Kok%NULL%2,                Kelvin Kai-Wang%To%NULL%2,                Hin%Chu%NULL%2,                Jin%Yang%NULL%2,                Fanfan%Xing%NULL%2,                Jieling%Liu%NULL%2,                Cyril Chik-Yan%Yip%NULL%3,                Rosana Wing-Shan%Poon%NULL%3,                Hoi-Wah%Tsoi%NULL%2,                Simon Kam-Fai%Lo%NULL%2,                Kwok-Hung%Chan%NULL%3,                Vincent Kwok-Man%Poon%NULL%2,                Wan-Mui%Chan%NULL%3,                Jonathan Daniel%Ip%NULL%3,                Jian-Piao%Cai%NULL%3,                Vincent Chi-Chung%Cheng%NULL%3,                Honglin%Chen%NULL%3,                Christopher Kim-Ming%Hui%NULL%2,                Kwok-Yung%Yuen%NULL%9]'
$ws.Range("E4").Value = '[Qun%Li%NULL%2,                Xuhua%Guan%NULL%1,                Peng%Wu%NULL%1,                Xiaoye%Wang%NULL%1,                Lei%Zhou%NULL%1,                Yeqing%Tong%NULL%1,                Ruiqi%Ren%NULL%1,                Kathy S.M.%Leung%NULL%1,                Eric H.Y.%Lau%NULL%1,                Jessica Y.%Wong%NULL%1,                Xuesen%Xing%NULL%1,                Nijuan%Xiang%NULL%1,                Yang%Wu%NULL%1,                Chao%Li%NULL%1,                Qi%Chen%NULL%1,                Dan%Li%NULL%1,                Tian%Liu%NULL%1,                Jing%Zhao%NULL%1,                Man%Liu%NULL%1,                Wenxiao%Tu%NULL%1,                Chuding%Chen%NULL%1,                Lianmei%Jin%NULL%1,                Rui%Yang%NULL%1,                Qi%Wang%NULL%1,                Suhua%Zhou%NULL%1,                Rui%Wang%NULL%1,                Hui%Liu%NULL%1,                Yinbo%Luo%NULL%1,                Yuan%Liu%NULL%1,                Ge%Shao%NULL%1,                Huan%Li%NULL%3,                Zhongfa%Tao%NULL%1,                Yang%Yang%NULL%3,                Zhiqiang%Deng%NULL%1,                Boxi%Liu%NULL%1,                Zhitao%Ma%NULL%1,                Yanping%Zhang%NULL%1,                Guoqing%Shi%NULL%1,                Tommy T.Y.%Lam%NULL%1,                Joseph T.%Wu%NULL%1,                George F.%Gao%NULL%1,                Benjamin J.%Cowling%NULL%1,                Bo%Yang%NULL%4,                Bo%Yang%NULL%0,                Gabriel M.%Leung%NULL%1,                Zijian%Feng%NULL%1]'
$ws.Range("E5").Value = '[Nanshan%Chen%NULL%9,                Min%Zhou%NULL%6,                Xuan%Dong%NULL%6,                Jieming%Qu%NULL%10,                Fengyun%Gong%NULL%6,                Yang%Han%NULL%5,                Yang%Qiu%NULL%5,                Jingli%Wang%NULL%6,                Ying%Liu%NULL%15,                Yuan%Wei%NULL%0,                Jia''an%Xia%NULL%11,                Ting%Yu%NULL%0,                Xinxin%Zhang%NULL%5,                Li%Zhang%NULL%0]'
$ws.Range("E6").Value = '[Fengxiang%Song%NULL%2,                Nannan%Shi%NULL%3,                Nannan%Shi%NULL%0,                Fei%Shan%NULL%2,                Zhiyong%Zhang%NULL%2,                Jie%Shen%NULL%1,                Hongzhou%Lu%NULL%6,                Yun%Ling%NULL%7,                Yebin%Jiang%NULL%2,                Yebin%Jiang%NULL%0,                Yuxin%Shi%shiyuxin@shphc.org.cn%5,                Yuxin%Shi%shiyuxin@shphc.org.cn%0]'
$ws.Range("E7").Value = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%11,            Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%5,            Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%5,            Subert-Salas%Lizandra%coreGivesNoEmail%5]'
$ws.Range("E8").Value = '[Kui%Liu%NULL%6,                Yuan-Yuan%Fang%NULL%4,                Yan%Deng%NULL%4,                Wei%Liu%NULL%6,                Mei-Fang%Wang%NULL%2,                Jing-Ping%Ma%NULL%2,                Wei%Xiao%NULL%2,                Ying-Nan%Wang%NULL%2,                Min-Hua%Zhong%NULL%2,                Cheng-Hong%Li%NULL%2,                Guang-Cai%Li%NULL%2,                Hui-Guo%Liu%NULL%6,                Xiu-Yuan%Hao%NULL%6,                Xiu-Yuan%Hao%NULL%0,                Pei-Fang%Wei%NULL%11]'
$ws.Range("E9").Value = '[De%Chang%xref no email%2,        Minggui%Lin%xref no email%2,        Lai%Wei%xref no email%2,        Lixin%Xie%xref no email%2,        Guangfa%Zhu%xref no email%2,        Charles S.%Dela Cruz%xref no email%2,        Lokesh%Sharma%xref no email%2]'
$ws.Range("E10").Value = '[ Yingxia%Liu%null%1,      Yang%Yang%null%1,      Cong%Zhang%null%1,      Fengming%Huang%null%1,      Fuxiang%Wang%null%1,      Jing%Yuan%null%2,      Zhaoqin%Wang%null%1,      Jinxiu%Li%null%1,      Jianming%Li%null%1,      Cheng%Feng%null%1,      Zheng%Zhang%null%2,      Lifei%Wang%null%1,      Ling%Peng%null%1,      Li%Chen%null%1,      Yuhao%Qin%null%1,      Dandan%Zhao%null%1,      Shuguang%Tan%null%1,      Lu%Yin%null%1,      Jun%Xu%null%1,      Congzhao%Zhou%null%1,      Chengyu%Jiang%null%1,      Lei%Liu%null%2]'
$ws.Range("E11").Value = '[Li-Li%Ren%NULL%3,                Ye-Ming%Wang%NULL%1,                Zhi-Qiang%Wu%NULL%1,                Zi-Chun%Xiang%NULL%1,                Li%Guo%NULL%0,                Teng%Xu%NULL%1,                Yong-Zhong%Jiang%NULL%1,                Yan%Xiong%NULL%4,                Yong-Jun%Li%NULL%1,                Xing-Wang%Li%NULL%1,                Hui%Li%NULL%0,                Guo-Hui%Fan%NULL%1,                Xiao-Ying%Gu%NULL%1,                Yan%Xiao%NULL%0,                Hong%Gao%NULL%0,                Jiu-Yang%Xu%NULL%1,                Fan%Yang%NULL%3,                Xin-Ming%Wang%NULL%1,                Chao%Wu%NULL%3,                Lan%Chen%NULL%1,                Yi-Wei%Liu%NULL%1,                Bo%Liu%NULL%2,                Jian%Yang%NULL%1,                Xiao-Rui%Wang%NULL%1,                Jie%Dong%NULL%1,                Li%Li%NULL%5,                Chao-Lin%Huang%NULL%1,                Jian-Ping%Zhao%NULL%1,                Yi%Hu%NULL%0,                Zhen-Shun%Cheng%NULL%1,                Lin-Lin%Liu%NULL%1,                Zhao-Hui%Qian%NULL%1,                Chuan%Qin%NULL%1,                Qi%Jin%NULL%0,                Bin%Cao%NULL%0,                Jian-Wei%Wang%NULL%1,                Xiu-Yuan%Hao%NULL%0,                Xiu-Yuan%Hao%NULL%0,                Pei-Fang%Wei%NULL%0]'
$ws.Range("E12").Value = '[Ning%Tang%NULL%4,                Dengju%Li%NULL%2,                Xiong%Wang%NULL%4,                Ziyong%Sun%NULL%8]'
$ws.Range("E13").Value = '[Akdis%Cezmi A%coreGivesNoEmail%9,            Cao%Yi-yuan%coreGivesNoEmail%4,            Dong%Xiang%coreGivesNoEmail%4,            Gao%Ya-dong%coreGivesNoEmail%4,            Yan%You-qin%coreGivesNoEmail%4,            Yang%Yi-bin%coreGivesNoEmail%4,            Yuan%Ya-dong%coreGivesNoEmail%4,            Zhang%Jin-jin%coreGivesNoEmail%4]'
$ws.Range("E14").Value = '[Xiao-Wei%Xu%NULL%7,                Xiao-Xin%Wu%NULL%3,                Xian-Gao%Jiang%NULL%3,                Kai-Jin%Xu%NULL%3,                Ling-Jun%Ying%NULL%3,                Chun-Lian%Ma%NULL%3,                Shi-Bo%Li%NULL%3,                Hua-Ying%Wang%NULL%3,                Sheng%Zhang%NULL%3,                Hai-Nv%Gao%NULL%3,                Ji-Fang%Sheng%NULL%3,                Hong-Liu%Cai%NULL%3,                Yun-Qing%Qiu%NULL%3,                Lan-Juan%Li%NULL%3]'
$ws.Range("E15").Value = '[Jiong%Wu%NULL%1,                Xiaojia%Wu%NULL%1,                Wenbing%Zeng%NULL%1,                Dajing%Guo%NULL%4,                Zheng%Fang%NULL%4,                Linli%Chen%NULL%4,                Huizhe%Huang%NULL%1,                Chuanming%Li%NULL%4]'
$ws.Range("E16").Value = '[Kin On%Kwok%kkokwok@cuhk.edu.hk%2,                Valerie%Wong%NULL%1,                Vivian Wan In%Wei%NULL%1,                Samuel Yeung Shan%Wong%NULL%1,                Julian Wei-Tze%Tang%julian.tang@uhl-tr.nhs.uk%1]'
$ws.Range("E17").Value = '[Heshui%Shi%NULL%2,                Xiaoyu%Han%NULL%1,                Nanchuan%Jiang%NULL%1,                Yukun%Cao%NULL%1,                Osamah%Alwalid%NULL%1,                Jin%Gu%NULL%1,                Yanqing%Fan%1024932023@qq.com%2,                Chuansheng%Zheng%hqzcsxh@sina.com%3]'
$ws.Range("E18").Value = '[Yu-Huan%Xu%NULL%2,                Jing-Hui%Dong%NULL%2,                Wei-Min%An%NULL%2,                Xiao-Yan%Lv%15001008285@139.com%2,                Xiao-Ping%Yin%yinxiaoping78@sina.com%2,                Jian-Zeng%Zhang%NULL%2,                Li%Dong%NULL%2,                Xi%Ma%NULL%2,                Hong-Jie%Zhang%NULL%2,                Bu-Lang%Gao%browngao@163.com%2]'
$ws.Range("E19").Value = '[Tao%Ai%NULL%2,                Zhenlu%Yang%NULL%2,                Zhenlu%Yang%NULL%0,                Hongyan%Hou%NULL%2,                Hongyan%Hou%NULL%0,                Chenao%Zhan%NULL%1,                Chong%Chen%NULL%1,                Wenzhi%Lv%NULL%2,                Wenzhi%Lv%NULL%0,                Qian%Tao%NULL%2,                Qian%Tao%NULL%0,                Ziyong%Sun%NULL%0,                Ziyong%Sun%NULL%0,                Liming%Xia%xialiming2017@outlook.com%1]'
$ws.Range("E20").Value = '[Wenjie%Yang%NULL%5,                Qiqi%Cao%NULL%3,                Le%Qin%NULL%3,                Xiaoyang%Wang%NULL%3,                Zenghui%Cheng%NULL%3,                Ashan%Pan%NULL%3,                Jianyi%Dai%NULL%3,                Qingfeng%Sun%NULL%3,                Fengquan%Zhao%NULL%3,                Jieming%Qu%NULL%0,                Fuhua%Yan%NULL%3]'
$ws.Range("E21").Value = '[Sijia%Tian%NULL%2,                Nan%Hu%NULL%2,                Jing%Lou%NULL%2,                Kun%Chen%NULL%2,                Xuqin%Kang%NULL%2,                Zhenjun%Xiang%NULL%2,                Hui%Chen%NULL%4,                Dali%Wang%NULL%2,                Ning%Liu%NULL%2,                Dong%Liu%NULL%2,                Gang%Chen%NULL%2,                Yongliang%Zhang%NULL%2,                Dou%Li%NULL%2,                Jianren%Li%NULL%2,                Huixin%Lian%NULL%2,                Shengmei%Niu%NULL%2,                Luxi%Zhang%NULL%2,                Jinjun%Zhang%NULL%2]'
$ws.Range("E22").Value = '[Yihui%Huang%NULL%3,                Mengqi%Tu%NULL%3,                Shipei%Wang%NULL%3,                Sichao%Chen%NULL%3,                Wei%Zhou%NULL%5,                Danyang%Chen%NULL%3,                Lin%Zhou%NULL%3,                Min%Wang%NULL%3,                Yan%Zhao%NULL%5,                Wen%Zeng%NULL%3,                Qi%Huang%NULL%3,                Hai''bo%Xu%NULL%3,                Zeming%Liu%NULL%3,                Liang%Guo%NULL%3]'
$ws.Range("E23").Value = '[Xi%Xu%NULL%4,                Chengcheng%Yu%NULL%2,                Jing%Qu%NULL%2,                Lieguang%Zhang%NULL%2,                Songfeng%Jiang%NULL%2,                Deyang%Huang%NULL%2,                Bihua%Chen%NULL%2,                Zhiping%Zhang%NULL%2,                Wanhua%Guan%NULL%2,                Zhoukun%Ling%NULL%2,                Rui%Jiang%NULL%2,                Tianli%Hu%NULL%2,                Yan%Ding%NULL%2,                Lin%Lin%NULL%2,                Qingxin%Gan%NULL%2,                Liangping%Luo%tluolp@jnu.edu.cn%2,                Xiaoping%Tang%xtang@21cn.com%2,                Jinxin%Liu%Liujx83710378@126.com%2]'
$ws.Range("E24").Value = '[Wei%Liu%NULL%8,                Zhao-Wu%Tao%NULL%2,                Lei%Wang%NULL%3,                Ming-Li%Yuan%NULL%4,                Kui%Liu%NULL%6,                Ling%Zhou%NULL%4,                Shuang%Wei%NULL%4,                Yan%Deng%NULL%0,                Jing%Liu%NULL%4,                Hui-Guo%Liu%NULL%0,                Ming%Yang%NULL%3,                Yi%Hu%NULL%0,                Pei-Fang%Wei%NULL%0,                Pei-Fang%Wei%NULL%0]'
$ws.Range("E25").Value = '[Kunhua%Li%NULL%3,                Jiong%Wu%NULL%3,                Faqi%Wu%NULL%3,                Dajing%Guo%NULL%0,                Linli%Chen%NULL%0,                Zheng%Fang%NULL%0,                Chuanming%Li%NULL%0]'
$ws.Range("E26").Value = '[Jian%Wu%NULL%2,                Jun%Liu%NULL%3,                Xinguo%Zhao%NULL%1,                Chengyuan%Liu%NULL%1,                Wei%Wang%NULL%6,                Dawei%Wang%NULL%1,                Wei%Xu%NULL%2,                Chunyu%Zhang%NULL%1,                Jiong%Yu%NULL%1,                Bin%Jiang%NULL%1,                Hongcui%Cao%hccao@zju.edu.cn%1,                Lanjuan%Li%NULL%1]'
$ws.Range("E27").Value = '[Jianlei%Cao%NULL%12,                Xiaorong%Hu%NULL%2,                Wenlin%Cheng%NULL%7,                Lei%Yu%NULL%7,                Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%7,                Qiang%Liu%liuqiang@irm-cams.ac.cn%9,                Qiang%Liu%liuqiang@irm-cams.ac.cn%0]'
$ws.Range("E28").Value = '[Qiurong%Ruan%NULL%6,                Kun%Yang%NULL%10,                Kun%Yang%NULL%0,                Wenxia%Wang%NULL%10,                Wenxia%Wang%NULL%0,                Lingyu%Jiang%NULL%10,                Lingyu%Jiang%NULL%0,                Jianxin%Song%songsingsjx@sina.com%10,                Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("E29").Value = '[Qiurong%Ruan%NULL%0,               Kun%Yang%NULL%3,               Kun%Yang%NULL%0,               Wenxia%Wang%NULL%2,               Wenxia%Wang%NULL%0,               Lingyu%Jiang%NULL%2,               Lingyu%Jiang%NULL%0,               Jianxin%Song%songsingsjx@sina.com%2,               Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("E30").Value = '[Fei%Xiao%NULL%1,                Meiwen%Tang%NULL%1,                Xiaobin%Zheng%NULL%1,                Ye%Liu%ly77219@163.com%1,                Xiaofeng%Li%zdwylxf@163.com%1,                Hong%Shan%shanhong@mail.sysu.edu.cn%13]'
$ws.Range("E31").Value = '[Ying%Xiong%NULL%2,                Dong%Sun%NULL%1,                Yao%Liu%NULL%1,                Yanqing%Fan%NULL%0,                Lingyun%Zhao%NULL%1,                Xiaoming%Li%NULL%1,                Wenzhen%Zhu%NULL%1]'
$ws.Range("E32").Value = '[Ying%Xiong%NULL%0,               Dong%Sun%NULL%1,               Yao%Liu%NULL%1,               Yanqing%Fan%NULL%1,               Lingyun%Zhao%NULL%1,               Xiaoming%Li%NULL%2,               Wenzhen%Zhu%NULL%1]'
$ws.Range("C33").Value = '"Epidemiologic Features and Clinical Course of Patients Infected With SARS-CoV-2 in Singapore"'
$ws.Range("D33").Value = 'Unknown Abstract'
$ws.Range("E33").Value = '[Barnaby Edward%Young%xref no email%0, Sean Wei Xiang%Ong%xref no email%1, Shirin%Kalimuddin%xref no email%1, Jenny G.%Low%xref no email%1, Seow Yen%Tan%xref no email%1, Jiashen%Loh%xref no email%1, Oon-Tek%Ng%xref no email%1, Kalisvar%Marimuthu%xref no email%1, Li Wei%Ang%xref no email%1, Tze Minn%Mak%xref no email%1, Sok Kiang%Lau%xref no email%1, Danielle E.%Anderson%xref no email%1, Kian Sing%Chan%xref no email%1, Thean Yen%Tan%xref no email%1, Tong Yong%Ng%xref no email%1, Lin%Cui%xref no email%1, Zubaidah%Said%xref no email%1, Lalitha%Kurupatham%xref no email%1, Mark I-Cheng%Chen%xref no email%1, Monica%Chan%xref no email%1, Shawn%Vasoo%xref no email%1, Lin-Fa%Wang%xref no email%1, Boon Huan%Tan%xref no email%1, Raymond Tzer Pin%Lin%xref no email%1, Vernon Jian Ming%Lee%xref no email%1, Yee-Sin%Leo%xref no email%1, David Chien%Lye%xref no email%1]'
$ws.Range("F33").Value = '10.1001/jama.2020.3204'
$ws.Range("G33").Value = 'CROSSREF'
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '2023-05-30'
$ws.Range("E34").Value = '[JingCheng%Zhang%zjc1983@126.com%0,               SaiBin%Wang%NULL%2,               SaiBin%Wang%NULL%0,               YaDong%Xue%NULL%1]'
$ws.Range("C35").Value = '"Coronavirus Disease 2019 (COVID-19): Role of Chest CT in Diagnosis and Management"'
$ws.Range("D35").Value = 'Unknown Abstract'
$ws.Range("E35").Value = '[Yan%Li%xref no email%1, Liming%Xia%xref no email%2]'
$ws.Range("F35").Value = '10.2214/AJR.20.22954'
$ws.Range("G35").Value = 'CROSSREF'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '2023-05-30'
$ws.Range("E36").Value = '[Yang%Yao%NULL%2,                Yao%Tian%NULL%2,                Yao%Tian%NULL%0,                Jing%Zhou%NULL%2,                Xuan%Ma%NULL%1,                Min%Yang%NULL%1,                ShengYu%Wang%NULL%1]'
$ws.Range("E37").Value = '[Yang%Yao%NULL%0,               Yao%Tian%NULL%2,               Yao%Tian%NULL%0,               Jing%Zhou%NULL%2,               Xuan%Ma%NULL%1,               Min%Yang%NULL%1,               ShengYu%Wang%NULL%1]'
$ws.Range("C38").Value = '"CT Features of Coronavirus Disease 2019 (COVID-19) Pneumonia in 62 Patients in Wuhan, China"'
$ws.Range("D38").Value = 'Unknown Abstract'
$ws.Range("E38").Value = '[Shuchang%Zhou%xref no email%1, Yujin%Wang%xref no email%1, Tingting%Zhu%xref no email%1, Liming%Xia%xref no email%0]'
$ws.Range("F38").Value = '10.2214/AJR.20.22975'
$ws.Range("G38").Value = 'CROSSREF'
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = '2023-05-30'
$ws.Range("E39").Value = '[Lei%Wang%NULL%1,                Yong-hua%Gao%NULL%1,                Li-Li%lou%NULL%1,                Guo-Jun%Zhang%NULL%1]'
$ws.Range("E40").Value = '[Simin%Zhang%NULL%1,                Huaqiao%Li%NULL%1,                Songtao%Huang%NULL%1,                Wei%You%NULL%1,                Huaiqiang%Sun%NULL%1]'
$ws.Range("E41").Value = '[Fei%Zhou%NULL%19,                Ting%Yu%NULL%0,                Ronghui%Du%NULL%9,                Guohui%Fan%NULL%0,                Ying%Liu%NULL%0,                Zhibo%Liu%NULL%10,                Jie%Xiang%NULL%16,                Yeming%Wang%NULL%0,                Bin%Song%NULL%10,                Xiaoying%Gu%NULL%0,                Lulu%Guan%NULL%9,                Yuan%Wei%NULL%0,                Hui%Li%NULL%0,                Xudong%Wu%NULL%10,                Jiuyang%Xu%NULL%0,                Shengjin%Tu%NULL%9,                Yi%Zhang%NULL%10,                Hua%Chen%NULL%10,                Bin%Cao%NULL%0]'
$ws.Range("E42").Value = '[Fei%Zhou%NULL%0,               Ting%Yu%NULL%16,               Ronghui%Du%NULL%7,               Guohui%Fan%NULL%10,               Ying%Liu%NULL%9,               Zhibo%Liu%NULL%7,               Jie%Xiang%NULL%10,               Yeming%Wang%NULL%10,               Bin%Song%NULL%7,               Xiaoying%Gu%NULL%10,               Lulu%Guan%NULL%7,               Yuan%Wei%NULL%11,               Hui%Li%NULL%14,               Xudong%Wu%NULL%7,               Jiuyang%Xu%NULL%10,               Shengjin%Tu%NULL%7,               Yi%Zhang%NULL%8,               Hua%Chen%NULL%7,               Bin%Cao%NULL%11]'
$ws.Range("C43").Value = '"COVID-19, Australia: Epidemiology Report 6: Reporting week ending 1900 AEDT 7 March 2020"'
$ws.Range("D43").Value = '"This is the sixth epidemiological report for coronavirus disease 2019 (COVID-19), reported in Australia as at 19:00 Australian Eastern Daylight Time [AEDT] 7 March 2020. It includes data on COVID-19 cases diagnosed in Australia, the international situation and a review of current evidence.
'
$ws.Range("E43").Value = '[Tracy%Tsang%xref no email%1]'
$ws.Range("F43").Value = '10.33321/cdi.2020.44.21'
$ws.Range("G43").Value = 'CROSSREF'
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = '2023-02-21'
$ws.Range("E44").Value = '[Jun%Chen%NULL%4,                Tangkai%Qi%NULL%2,                Li%Liu%NULL%2,                Yun%Ling%NULL%0,                Zhiping%Qian%NULL%2,                Tao%Li%NULL%3,                Feng%Li%NULL%2,                Qingnian%Xu%NULL%2,                Yuyi%Zhang%NULL%2,                Shuibao%Xu%NULL%2,                Zhigang%Song%NULL%2,                Yigang%Zeng%NULL%2,                Yinzhong%Shen%NULL%2,                Yuxin%Shi%NULL%0,                Tongyu%Zhu%NULL%3,                Hongzhou%Lu%NULL%0]'
$ws.Range("E45").Value = '[Dahai%Zhao%NULL%2,                Feifei%Yao%NULL%1,                Lijie%Wang%NULL%1,                Ling%Zheng%NULL%1,                Yongjun%Gao%NULL%1,                Jun%Ye%NULL%1,                Feng%Guo%NULL%1,                Hui%Zhao%NULL%1,                Rongbao%Gao%rongbaogao@163.com%1]'
$ws.Range("E46").Value = '[Qing%Chen%NULL%2,                Bin%Quan%NULL%1,                Xiaoning%Li%NULL%2,                Guangjian%Gao%NULL%1,                Wenqiang%Zheng%NULL%1,                Jun%Zhang%NULL%1,                Zhiyun%Zhang%NULL%1,                Chunsheng%Liu%NULL%1,                Li%Li%NULL%0,                Chenglin%Wang%NULL%1,                Guihua%Zhang%NULL%1,                Jiajia%Li%NULL%1,                Yunhai%Dai%NULL%1,                Jianghua%Yang%yjhpath@163.com%1,                Wenzheng%Han%11418166@zju.edu.cn%1]'
$ws.Range("E47").Value = '[Qing%Chen%NULL%0,               Bin%Quan%NULL%1,               Xiaoning%Li%NULL%2,               Guangjian%Gao%NULL%1,               Wenqiang%Zheng%NULL%1,               Jun%Zhang%NULL%5,               Zhiyun%Zhang%NULL%1,               Chunsheng%Liu%NULL%1,               Li%Li%NULL%4,               Chenglin%Wang%NULL%1,               Guihua%Zhang%NULL%1,               Jiajia%Li%NULL%1,               Yunhai%Dai%NULL%1,               Jianghua%Yang%yjhpath@163.com%1,               Wenzheng%Han%11418166@zju.edu.cn%1]'
$ws.Range("C48").Value = '"Risk Factors Associated With Acute Respiratory Distress Syndrome and Death in Patients With Coronavirus Disease 2019 Pneumonia in Wuhan, China"'
$ws.Range("D48").Value = 'Unknown Abstract'
$ws.Range("E48").Value = '[Chaomin%Wu%xref no email%0, Xiaoyan%Chen%xref no email%1, Yanping%Cai%xref no email%1, Jia\u2019an%Xia%xref no email%1, Xing%Zhou%xref no email%1, Sha%Xu%xref no email%1, Hanping%Huang%xref no email%1, Li%Zhang%xref no email%0, Xia%Zhou%xref no email%1, Chunling%Du%xref no email%1, Yuye%Zhang%xref no email%1, Juan%Song%xref no email%1, Sijiao%Wang%xref no email%1, Yencheng%Chao%xref no email%1, Zeyong%Yang%xref no email%1, Jie%Xu%xref no email%1, Xin%Zhou%xref no email%1, Dechang%Chen%xref no email%1, Weining%Xiong%xref no email%1, Lei%Xu%xref no email%1, Feng%Zhou%xref no email%1, Jinjun%Jiang%xref no email%1, Chunxue%Bai%xref no email%1, Junhua%Zheng%xref no email%1, Yuanlin%Song%xref no email%1]'
$ws.Range("F48").Value = '10.1001/jamainternmed.2020.0994'
$ws.Range("G48").Value = 'CROSSREF'
$ws.Range("H48").NumberFormat = "@"
$ws.Range("H48").Value = '2023-05-30'
$ws.Range("E49").Value = '[Wanbo%Zhu%NULL%0,               Kai%Xie%NULL%2,               Kai%Xie%NULL%0,               Hui%Lu%NULL%1,               Lei%Xu%bayinhexl@126.com%1,               Shusheng%Zhou%zhouss108@163.com%1,               Shiyuan%Fang%fangshiyuan2008@126.com%1]'
$ws.Range("C50").Value = '"Clinical Features and Chest CT Manifestations of Coronavirus Disease 2019 (COVID-19) in a Single-Center Study in Shanghai, China"'
$ws.Range("D50").Value = 'Unknown Abstract'
$ws.Range("E50").Value = '[Zenghui%Cheng%xref no email%0, Yong%Lu%xref no email%1, Qiqi%Cao%xref no email%1, Le%Qin%xref no email%1, Zilai%Pan%xref no email%1, Fuhua%Yan%xref no email%1, Wenjie%Yang%xref no email%0]'
$ws.Range("F50").Value = '10.2214/AJR.20.22959'
$ws.Range("G50").Value = 'CROSSREF'
$ws.Range("H50").NumberFormat = "@"
$ws.Range("H50").Value = '2023-05-09'
$ws.Range("E51").Value = '[Zhongliang%Wang%NULL%6,                Bohan%Yang%NULL%3,                Qianwen%Li%NULL%3,                Lu%Wen%NULL%3,                Ruiguang%Zhang%zrg27@163.com%3]'
$ws.Range("E52").Value = '[Yanrong%Wang%123rong@sohu.com%2,                Yingxia%Liu%NULL%1,                Lei%Liu%NULL%15,                Xianfeng%Wang%NULL%1,                Nijuan%Luo%NULL%1,                Ling%Li%NULL%1]'
$ws.Range("E53").Value = '[Yanrong%Wang%123rong@sohu.com%0,               Yingxia%Liu%NULL%2,               Lei%Liu%NULL%9,               Xianfeng%Wang%NULL%1,               Nijuan%Luo%NULL%1,               Ling%Li%NULL%1]'
$ws.Range("C54").Value = '"Early Clinical and CT Manifestations of Coronavirus Disease 2019 (COVID-19) Pneumonia"'
$ws.Range("D54").Value = 'Unknown Abstract'
$ws.Range("E54").Value = '[Rui%Han%xref no email%1, Lu%Huang%xref no email%1, Hong%Jiang%xref no email%1, Jin%Dong%xref no email%1, Hongfen%Peng%xref no email%1, Dongyou%Zhang%xref no email%1]'
$ws.Range("F54").Value = '10.2214/AJR.20.22961'
$ws.Range("G54").Value = 'CROSSREF'
$ws.Range("H54").NumberFormat = "@"
$ws.Range("H54").Value = '2023-05-10'
$ws.Range("E55").Value = '[Yong%Gao%NULL%4,                Tuantuan%Li%NULL%2,                Mingfeng%Han%NULL%2,                Xiuyong%Li%NULL%2,                Dong%Wu%NULL%2,                Yuanhong%Xu%NULL%4,                Yulin%Zhu%NULL%2,                Yan%Liu%NULL%3,                Xiaowu%Wang%wangxiaowu19880218@126.com%2,                Linding%Wang%wanglinding@ahmu.edu.cn%2]'
$ws.Range("E56").Value = '[Yong%Gao%NULL%0,               Tuantuan%Li%NULL%2,               Mingfeng%Han%NULL%6,               Xiuyong%Li%NULL%2,               Dong%Wu%NULL%2,               Yuanhong%Xu%NULL%2,               Yulin%Zhu%NULL%2,               Yan%Liu%NULL%2,               Xiaowu%Wang%wangxiaowu19880218@126.com%2,               Linding%Wang%wanglinding@ahmu.edu.cn%2]'
